# Generate Report for Handback
# -----------------------------------------------------------------------
# This script reflects the moment a localization handback has been
# processed: the status flips from "Ready for handoff" to
# "Handed back: in sync with en-US" everywhere it is shown, and the
# per-language detail sheets (zh-cn / de-de) get their "Latest Target
# File" / "Latest Handback File" / "Latest Handback DateTime" columns
# populated with the freshly generated report data.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"
$mdFileName = "9d613ea1-a056-42a2-82ce-6bf0c780f2a4.md"
$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f4463efc5e0a5ee2e4066ea61f044e973ca8aa89/e2e/9d613ea1-a056-42a2-82ce-6bf0c780f2a4.md"

# --- Overview sheet: status columns for zh-cn (E2) and de-de (F2) ------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText

# widen the (now longer) status columns to fit the new text
$overview.Columns.Item(5).ColumnWidth = 29.16
$overview.Columns.Item(6).ColumnWidth = 29.16

# --- zh-cn detail sheet --------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $statusText

# Latest Target File: link back to the source markdown
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdFileName)

# Latest Handback File / DateTime
$zhcn.Range("J2").Value = "9d613ea1-a056-42a2-82ce-6bf0c780f2a4.812b40031a2078c6f0852f659780606da3ad265e.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-22 10:39:31"

$zhcn.Columns.Item(3).ColumnWidth = 29.16
$zhcn.Columns.Item(9).ColumnWidth = 39.16
$zhcn.Columns.Item(10).ColumnWidth = 39.16

# --- de-de detail sheet --------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $statusText

# Latest Target File: link back to the source markdown
$dede.Hyperlinks.Add($dede.Range("I2"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdFileName)

# Latest Handback File / DateTime
$dede.Range("J2").Value = "9d613ea1-a056-42a2-82ce-6bf0c780f2a4.812b40031a2078c6f0852f659780606da3ad265e.de-de.xlf"
$dede.Range("K2").Value = "2016-08-22 10:39:38"

$dede.Columns.Item(3).ColumnWidth = 29.16
$dede.Columns.Item(9).ColumnWidth = 39.16
$dede.Columns.Item(10).ColumnWidth = 39.16

Write-Host "Handback report generated."
